$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the URL shown/stored in D2 (shared string text); the cell keeps its
# existing hyperlink relationship (still pointing at http://192.168.168.111/
# in the rels part per the diff, only the displayed text changes).
$ws.Range("D2").Value = "http://192.168.168.107/"

# Update the sheet view: select P2, then scroll so column H is the
# left-most visible column (order matters - selecting re-centers the
# viewport, so the scroll position must be (re-)applied afterwards).
$ws.Range("P2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
